$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.906.40'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '1.889.38'
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '0.8202'
$ws.Range("E5").Value = '  +6.53%  '
$ws.Range("D6").Value = '241.35'
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '0.3213'
$ws.Range("E8").Value = '  +5.21%  '
$ws.Range("D9").Value = '26.37'
$ws.Range("E9").Value = '  +3.58%  '
$ws.Range("D10").Value = '0.07009'
$ws.Range("E10").Value = '  +2.31%  '
$ws.Range("D11").Value = '0.08027'
$ws.Range("E11").Value = '  +0.56%  '
$ws.Range("D12").Value = '0.7440'
$ws.Range("E12").Value = '  +0.92%  '
$ws.Range("D13").Value = '1.894.19'
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D14").Value = '5.185'
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").Value = '91.96'
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Value = '29.907.93'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").Value = '13.97'
$ws.Range("E17").Value = '  +1.56%  '
$ws.Range("D18").Value = '5.871'
$ws.Range("E18").Value = '  -0.41%  '
$ws.Range("D19").Value = '244.09'
$ws.Range("E19").Value = '  -0.89%  '
$ws.Range("E20").Value = '  +0.57%  '
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").Value = '2.143.18'
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").Value = '6.899'
$ws.Range("E24").Value = '  -0.66%  '
$ws.Range("D25").Value = '0.1551'
$ws.Range("E25").Value = '  +20.70%  '
$ws.Range("E26").Value = '  -0.47%  '
$ws.Range("D27").Value = '9.163'
$ws.Range("E27").Value = '  -1.13%  '
$ws.Range("D28").Value = '18.79'
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("D29").Value = '2.075'
$ws.Range("E29").Value = '  +2.37%  '
$ws.Range("E30").Value = '  -1.86%  '
$ws.Range("E31").Value = '  +0.65%  '
$ws.Range("D32").Value = '4.263'
$ws.Range("E32").Value = '  -0.15%  '
$ws.Range("D33").Value = '0.05621'
$ws.Range("E33").Value = '  +7.09%  '
$ws.Range("D34").Value = '4.062'
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("E35").Value = '  +1.83%  '
$ws.Range("D36").Value = '0.7281'
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("D37").Value = '2.723'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("E38").Value = '  +0.22%  '
$ws.Range("D39").Value = '2.772'
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("D40").Value = '0.4406'
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").Value = '71.62'
$ws.Range("E41").Value = '  -0.46%  '
$ws.Range("D42").Value = '5.951'
$ws.Range("E42").Value = '  -3.83%  '
$ws.Range("D43").Value = '0.8441'
$ws.Range("E43").Value = '  +0.91%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D46").Value = '100.73'
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("D47").Value = '7.557'
$ws.Range("E47").Value = '  -0.46%  '
$ws.Range("D48").Value = '9.664'
$ws.Range("E48").Value = '  -0.96%  '
$ws.Range("D49").Value = '986.82'
$ws.Range("E49").Value = '  +7.73%  '
$ws.Range("D50").Value = '2.041.89'
$ws.Range("E50").Value = '  -0.69%  '
$ws.Range("D51").Value = '35.95'
$ws.Range("E51").Value = '  -0.50%  '
